# Updates the "剩余" (remaining) column E, and for a handful of rows that
# rolled over (remaining hit 1 and got renewed) the "开始时间" (start date)
# column F as well. Mirrors a daily countdown refresh of the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: E = 5 (剩余/remaining days), F = 6 (开始时间/start date)
$colE = 5
$colF = 6

# row, new E value, new F value (blank = leave F unchanged)
$updates = @"
2,8,
3,8,
4,8,
5,8,
6,8,
7,8,
8,8,
9,8,
10,1,
11,8,
12,8,
13,8,
14,8,
15,8,
16,2,
17,8,
18,1,
19,1,
20,1,
21,1,
22,8,
23,8,
24,8,
25,8,
26,8,
27,2,
28,1,
29,1,
30,1,
31,1,
32,1,
33,1,
34,1,
35,1,
37,1,
38,1,
39,1,
40,1,
41,1,
42,1,
43,8,
44,1,
45,8,
46,1,
47,1,
48,1,
49,2,
50,6,
51,6,
52,6,
53,6,
54,6,
55,6,
56,6,
57,6,
58,10,20260215
59,10,20260215
60,10,20260215
61,2,
62,10,20260215
63,10,20260215
64,10,20260215
65,1,
66,1,
67,1,
68,1,
69,1,
70,2,
71,2,
72,2,
73,2,
74,2,
75,2,
76,2,
77,5,
78,5,
79,5,
80,5,
81,5,
82,5,
83,5,
84,5,
85,5,
86,5,
87,1,
88,1,
89,1,
90,1,
91,8,
92,1,
93,5,
94,4,
95,4,
96,2,
97,2,
98,2,
99,2,
"@

$rows = $updates -split "`n" | Where-Object { $_.Trim() -ne "" }

foreach ($line in $rows) {
    $parts = $line.Trim() -split ","
    $rowNum = [int]$parts[0]
    $newE = [int]$parts[1]

    $ws.Cells.Item($rowNum, $colE).Value = $newE

    if ($parts.Length -ge 3 -and $parts[2].Trim() -ne "") {
        $newF = [int]$parts[2]
        $ws.Cells.Item($rowNum, $colF).Value = $newF
    }
}
